# ============================================================
# Edit: Truncate overly long population/sample names
# 1. Add "Truncated Names" lookup sheet (Long Name -> Short Name)
# 2. Replace truncated Code values in pop_names!A with short codes
# ============================================================

$wb = $excel.ActiveWorkbook

# --- Step 1: update pop_names Code column (A) with short codes ---
$popNames = $wb.Worksheets.Item("pop_names")
$popNames.Range("A3").Value = "BAS2"
$popNames.Range("A5").Value = "BEA"
$popNames.Range("A6").Value = "BOX"
$popNames.Range("A11").Value = "AM"
$popNames.Range("A17").Value = "EUR"
$popNames.Range("A18").Value = "FS"
$popNames.Range("A19").Value = "GSD"
$popNames.Range("A20").Value = "SAM"
$popNames.Range("A21").Value = "SH"
$popNames.Range("A22").Value = "AED"
$popNames.Range("A23").Value = "APBT"
$popNames.Range("A24").Value = "AST"
$popNames.Range("A25").Value = "CD"
$popNames.Range("A26").Value = "CLD"
$popNames.Range("A27").Value = "CBR"
$popNames.Range("A28").Value = "CHI"
$popNames.Range("A31").Value = "NEW"
$popNames.Range("A32").Value = "NSDTR"
$popNames.Range("A33").Value = "PIO"
$popNames.Range("A34").Value = "VDB"
$popNames.Range("A35").Value = "VDB2"
$popNames.Range("A36").Value = "VDC"
$popNames.Range("A37").Value = "VDCR"
$popNames.Range("A38").Value = "VDDR"
$popNames.Range("A39").Value = "VDH"
$popNames.Range("A40").Value = "VDP"
$popNames.Range("A41").Value = "VDPA"
$popNames.Range("A42").Value = "VDPC"
$popNames.Range("A43").Value = "VDPI"
$popNames.Range("A44").Value = "VDPL"
$popNames.Range("A45").Value = "VDPP"
$popNames.Range("A46").Value = "VDPR"
$popNames.Range("A47").Value = "VDUA"
$popNames.Range("A48").Value = "XOL"
$popNames.Range("A51").Value = "VDIC"
$popNames.Range("A52").Value = "VDID"
$popNames.Range("A53").Value = "VDIH"
$popNames.Range("A54").Value = "VDIM"
$popNames.Range("A55").Value = "VDIO"
$popNames.Range("A56").Value = "CSP"
$popNames.Range("A57").Value = "CC"
$popNames.Range("A62").Value = "NGSD"
$popNames.Range("A63").Value = "VDIB"
$popNames.Range("A64").Value = "VDIJ"
$popNames.Range("A65").Value = "VDPNGEH"
$popNames.Range("A66").Value = "VDPNGPM"
$popNames.Range("A67").Value = "VDVCB"
$popNames.Range("A68").Value = "VDVHG"
$popNames.Range("A69").Value = "VDVLS"
$popNames.Range("A70").Value = "VDVLC"
$popNames.Range("A79").Value = "TAI"

# Autofit column A to reflect the now-short Code values (adds <col> entry)
$popNames.Columns.Item(1).AutoFit() | Out-Null

# --- Step 2: create the new "Truncated Names" worksheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$tn.Name = "Truncated Names"

# Populate A1:B52 with the long-name/short-name lookup table (already alphabetised)
$tnData = New-Object 'object[,]' 52,2
$tnData[0,0] = "Long Name"
$tnData[0,1] = "Short Name"
$tnData[1,0] = "Alaskan_Malamute "
$tnData[1,1] = "AM"
$tnData[2,0] = "American_Eskimo_Dog "
$tnData[2,1] = "AED"
$tnData[3,0] = "American_Pit_Bull_Terrier "
$tnData[3,1] = "APBT"
$tnData[4,0] = "American_Staffordshire_Terrier "
$tnData[4,1] = "AST"
$tnData[5,0] = "Basenji "
$tnData[5,1] = "BAS2"
$tnData[6,0] = "Beagle "
$tnData[6,1] = "BEA"
$tnData[7,0] = "Boxer "
$tnData[7,1] = "BOX"
$tnData[8,0] = "Carolina_Dog "
$tnData[8,1] = "CD"
$tnData[9,0] = "Catahoula_Leopard_Dog "
$tnData[9,1] = "CLD"
$tnData[10,0] = "Chesapeake_Bay_Retriever "
$tnData[10,1] = "CBR"
$tnData[11,0] = "Chihuahua "
$tnData[11,1] = "CHI"
$tnData[12,0] = "Chinese_Shar-pei "
$tnData[12,1] = "CSP"
$tnData[13,0] = "Chow_Chow "
$tnData[13,1] = "CC"
$tnData[14,0] = "Eurasier "
$tnData[14,1] = "EUR"
$tnData[15,0] = "Finnish_Spitz "
$tnData[15,1] = "FS"
$tnData[16,0] = "Greenland_Sledge_Dog "
$tnData[16,1] = "GSD"
$tnData[17,0] = "New_Guinea_Singing_Dog "
$tnData[17,1] = "NGSD"
$tnData[18,0] = "Newfoundland "
$tnData[18,1] = "NEW"
$tnData[19,0] = "Nova_Scotia_Duck_Tolling_Retriever "
$tnData[19,1] = "NSDTR"
$tnData[20,0] = "Peruvian_Inca_Orchid "
$tnData[20,1] = "PIO"
$tnData[21,0] = "Samoyed "
$tnData[21,1] = "SAM"
$tnData[22,0] = "Siberian_Husky "
$tnData[22,1] = "SH"
$tnData[23,0] = "Taimyr "
$tnData[23,1] = "TAI"
$tnData[24,0] = "Village_Dog_Belize "
$tnData[24,1] = "VDB"
$tnData[25,0] = "Village_Dog_Brazil "
$tnData[25,1] = "VDB2"
$tnData[26,0] = "Village_Dog_Colombia "
$tnData[26,1] = "VDC"
$tnData[27,0] = "Village_Dog_Costa_Rica "
$tnData[27,1] = "VDCR"
$tnData[28,0] = "Village_Dog_Dominican_Republic "
$tnData[28,1] = "VDDR"
$tnData[29,0] = "Village_Dog_Honduras "
$tnData[29,1] = "VDH"
$tnData[30,0] = "Village_Dog_India-Chennai "
$tnData[30,1] = "VDIC"
$tnData[31,0] = "Village_Dog_India-Dehli "
$tnData[31,1] = "VDID"
$tnData[32,0] = "Village_Dog_India-Hazaribagh "
$tnData[32,1] = "VDIH"
$tnData[33,0] = "Village_Dog_India-Mumbai "
$tnData[33,1] = "VDIM"
$tnData[34,0] = "Village_Dog_India-Orissa "
$tnData[34,1] = "VDIO"
$tnData[35,0] = "Village_Dog_Indonesia-Borneo "
$tnData[35,1] = "VDIB"
$tnData[36,0] = "Village_Dog_Indonesia-Jakarta "
$tnData[36,1] = "VDIJ"
$tnData[37,0] = "Village_Dog_Panama "
$tnData[37,1] = "VDP"
$tnData[38,0] = "Village_Dog_Papua_New_Guinea-East_Highlands_ "
$tnData[38,1] = "VDPNGEH"
$tnData[39,0] = "Village_Dog_Papua_New_Guinea-Port_Moresby "
$tnData[39,1] = "VDPNGPM"
$tnData[40,0] = "Village_Dog_Peru-Arequipa "
$tnData[40,1] = "VDPA"
$tnData[41,0] = "Village_Dog_Peru-Cusco "
$tnData[41,1] = "VDPC"
$tnData[42,0] = "Village_Dog_Peru-Ica "
$tnData[42,1] = "VDPI"
$tnData[43,0] = "Village_Dog_Peru-Loreto "
$tnData[43,1] = "VDPL"
$tnData[44,0] = "Village_Dog_Peru-Puno "
$tnData[44,1] = "VDPP"
$tnData[45,0] = "Village_Dog_Puerto_Rico "
$tnData[45,1] = "VDPR"
$tnData[46,0] = "Village_Dog_US-Alaska "
$tnData[46,1] = "VDUA"
$tnData[47,0] = "Village_Dog_Vietnam-Cao_Bang "
$tnData[47,1] = "VDVCB"
$tnData[48,0] = "Village_Dog_Vietnam-Ha_Giang "
$tnData[48,1] = "VDVHG"
$tnData[49,0] = "Village_Dog_Vietnam-Lang_Son "
$tnData[49,1] = "VDVLS"
$tnData[50,0] = "Village_Dog_Vietnam-Lao_Cai "
$tnData[50,1] = "VDVLC"
$tnData[51,0] = "Xoloitzcuintli "
$tnData[51,1] = "XOL"
$tn.Range("A1:B52").Value = $tnData

# Column widths to fit the long-name / short-name text
$tn.Columns.Item(1).AutoFit() | Out-Null
$tn.Columns.Item(2).AutoFit() | Out-Null

# Apply AutoFilter over the table and record the sort state (data already sorted by Long Name)
$tnRange = $tn.Range("A1:B52")
$tnRange.AutoFilter()
$tn.Sort.SortFields.Clear()
$tn.Sort.SortFields.Add($tn.Range("A1:A52"))
$tn.Sort.SetRange($tnRange)
$tn.Sort.Header = 1
$tn.Sort.Apply()

# Register the hidden _FilterDatabase defined name for the new sheet (mirrors pop_names)
$tnFilterName = $tn.Names.Add("_xlnm._FilterDatabase", "='Truncated Names'!`$A`$1:`$B`$52", $true)
$tnFilterName.Visible = $false

# Set the active selections to match the target workbook state
$tn.Range("G5").Select()
$popNames.Activate()
$popNames.Range("B6").Select()
